# Regenerate the "K" column (strikeouts) for each start, row-by-row,
# replacing the previously-written "Strike#" values that were stored in
# column G. Rows 2-52 correspond to data row indices 0-50.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2,0,1,2,0,1,1,1,1,2,0,0,1,0,0,0,0,1,0,1,2,2,2,1,0,1,0,0,1,0,0,2,1,3,2,1,2,4,0,3,0,3,0,2,1,1,2,3,2,1,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
